# AWS Qwiklabs.docx edit script
# Applies 4 changes described by the commit diff:
#  1. "Search for the lab by name ..." heading -> italic + color F17E3A (was bold/red/underline)
#  2. "Lab 4: ..." paragraph -> trailing space + bold/red/underline "DONE" run
#  3. "Lab 5: ..." paragraph -> two trailing spaces + bold/red/underline "DONE" run
#  4. "Lab 6: ..." paragraph -> split into "Lab 6: " / bookmark / "Introduction to Amazon DynamoDB"

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Search for the lab by name - they keep changing the URLs"
# Re-colour from CC0000 to F17E3A and make the whole paragraph (incl. mark) italic.
# ---------------------------------------------------------------------------
$headingRange = $d.Content.Duplicate
$headingRange.Find.ClearFormatting()
$headingRange.Find.Text = "Search for the lab by name"
$headingFound = $headingRange.Find.Execute()
if ($headingFound) {
    $headingPara = $headingRange.Paragraphs(1)
    $headingParaRange = $headingPara.Range
    $headingParaRange.Font.Italic = $true
    $headingParaRange.Font.ItalicBi = $true
    $headingParaRange.Font.Color = 3833585   # wdColor equivalent of F17E3A (RGB)
}

# ---------------------------------------------------------------------------
# Change 2: "Lab 4: Introduction to Amazon Simple Storage Service"
# Add a single trailing space, then a bold/red/underlined "DONE" run.
# ---------------------------------------------------------------------------
$lab4Range = $d.Content.Duplicate
$lab4Range.Find.ClearFormatting()
$lab4Range.Find.Text = "Lab 4: Introduction to Amazon Simple Storage Service"
$lab4Found = $lab4Range.Find.Execute()
if ($lab4Found) {
    $lab4Range.Collapse(0)
    $lab4Range.InsertAfter(" ")
}

$lab4DoneRange = $d.Content.Duplicate
$lab4DoneRange.Find.ClearFormatting()
$lab4DoneRange.Find.Text = "Lab 4: Introduction to Amazon Simple Storage Service "
$lab4DoneFound = $lab4DoneRange.Find.Execute()
if ($lab4DoneFound) {
    $lab4DoneRange.Collapse(0)
    $lab4DoneRange.InsertAfter("DONE")
}

$lab4MarkerRange = $d.Content.Duplicate
$lab4MarkerRange.Find.ClearFormatting()
$lab4MarkerRange.Find.Text = "Lab 4: Introduction to Amazon Simple Storage Service DONE"
$lab4MarkerFound = $lab4MarkerRange.Find.Execute()
if ($lab4MarkerFound) {
    $lab4MarkerRange.Collapse(0)
    $lab4MarkerRange.MoveStart(1, -4)
    $lab4MarkerRange.Font.Bold = $true
    $lab4MarkerRange.Font.Color = 204        # wdColor equivalent of CC0000 (RGB)
    $lab4MarkerRange.Font.Underline = 1      # wdUnderlineSingle
}

# ---------------------------------------------------------------------------
# Change 3: "Lab 5: Introduction to Amazon Relational Database Service (RDS) (Linux)"
# Add two trailing spaces, then a bold/red/underlined "DONE" run.
# ---------------------------------------------------------------------------
$lab5Range = $d.Content.Duplicate
$lab5Range.Find.ClearFormatting()
$lab5Range.Find.Text = "Lab 5: Introduction to Amazon Relational Database Service (RDS) (Linux)"
$lab5Found = $lab5Range.Find.Execute()
if ($lab5Found) {
    $lab5Range.Collapse(0)
    $lab5Range.InsertAfter("  ")
}

$lab5DoneRange = $d.Content.Duplicate
$lab5DoneRange.Find.ClearFormatting()
$lab5DoneRange.Find.Text = "Lab 5: Introduction to Amazon Relational Database Service (RDS) (Linux)  "
$lab5DoneFound = $lab5DoneRange.Find.Execute()
if ($lab5DoneFound) {
    $lab5DoneRange.Collapse(0)
    $lab5DoneRange.InsertAfter("DONE")
}

$lab5MarkerRange = $d.Content.Duplicate
$lab5MarkerRange.Find.ClearFormatting()
$lab5MarkerRange.Find.Text = "Lab 5: Introduction to Amazon Relational Database Service (RDS) (Linux)  DONE"
$lab5MarkerFound = $lab5MarkerRange.Find.Execute()
if ($lab5MarkerFound) {
    $lab5MarkerRange.Collapse(0)
    $lab5MarkerRange.MoveStart(1, -4)
    $lab5MarkerRange.Font.Bold = $true
    $lab5MarkerRange.Font.Color = 204        # wdColor equivalent of CC0000 (RGB)
    $lab5MarkerRange.Font.Underline = 1      # wdUnderlineSingle
}

# ---------------------------------------------------------------------------
# Change 4: "Lab 6: Introduction to Amazon DynamoDB" -> split with a bookmark
# between "Lab 6: " and "Introduction to Amazon DynamoDB".
# ---------------------------------------------------------------------------
$splitRange = $d.Content.Duplicate
$splitRange.Find.ClearFormatting()
$splitRange.Find.Text = "Lab 6: "
$splitFound = $splitRange.Find.Execute()
if ($splitFound) {
    $splitRange.Collapse(0)
    $d.Bookmarks.Add("__DdeLink__100_4046424704", $splitRange)
}
